$d = $word.ActiveDocument

$replacements = @(
    @("2024-08-25 Sunday", "2024-08-26 Monday"),
    @("37÷6=", "31÷6="),
    @("67÷9=", "57÷2="),
    @("10÷5=", "43÷5="),
    @("60÷6=", "93÷7="),
    @("66÷2=", "14÷6="),
    @("30÷6=", "97÷3="),
    @("88÷6=", "60÷7="),
    @("69÷7=", "44÷2="),
    @("15÷9=", "39÷6="),
    @("75÷4=", "42÷6="),
    @("98÷3=", "32÷5="),
    @("24÷3=", "95÷4="),
    @("79÷3=", "33÷8="),
    @("81÷9=", "13÷2="),
    @("43÷2=", "12÷6="),
    @("99÷2=", "47÷2="),
    @("51÷4=", "71÷8="),
    @("77÷2=", "65÷9="),
    @("43÷7=", "38÷5="),
    @("93÷9=", "21÷7="),
    @("75÷7=", "65÷2="),
    @("97÷7=", "20÷9="),
    @("24÷8=", "23÷6="),
    @("83÷2=", "59÷6="),
    @("25÷3=", "90÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
